$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 21075
$ws.Range("J21").Value = 24150
$ws.Range("L21").Value = 24150
$ws.Range("N21").Value = -25086
$ws.Range("H23").Value = 21075
$ws.Range("J23").Value = 24150
$ws.Range("L23").Value = 24150
$ws.Range("N23").Value = -24618
$ws.Range("H34").Value = 9206.666999999999
$ws.Range("I34").Value = 810
$ws.Range("J34").Value = 26000
$ws.Range("K34").Value = 810
$ws.Range("L34").Value = 26000
$ws.Range("M34").Value = -607
$ws.Range("N34").Value = -26406
$ws.Range("H36").Value = 9206.666999999999
$ws.Range("I36").Value = 810
$ws.Range("J36").Value = 26000
$ws.Range("K36").Value = 810
$ws.Range("L36").Value = 26000
$ws.Range("M36").Value = -95
$ws.Range("N36").Value = -27430
$ws.Range("H47").Value = 9500
$ws.Range("I47").Value = 9500
$ws.Range("K47").Value = 9500
$ws.Range("M47").Value = -8528
$ws.Range("H76").Value = 6568.136
$ws.Range("I76").Value = 5158.3335
$ws.Range("J76").Value = 8259.9
$ws.Range("K76").Value = 5158.3335
$ws.Range("L76").Value = 8259.9
$ws.Range("M76").Value = -4843.3335
$ws.Range("N76").Value = -8889.9
$ws.Range("H79").Value = 6568.136
$ws.Range("I79").Value = 5158.3335
$ws.Range("J79").Value = 8259.9
$ws.Range("K79").Value = 5158.3335
$ws.Range("L79").Value = 8259.9
$ws.Range("M79").Value = -4066.3335
$ws.Range("N79").Value = -10443.9
$ws.Range("H141").Value = 1588.0646
$ws.Range("I141").Value = 1086.2963
$ws.Range("K141").Value = 3258.8889
$ws.Range("M141").Value = 1921.1111

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8672.293
$ws.Range("I32").Value = 8008.4688
$ws.Range("J32").Value = 12534.546
$ws.Range("K32").Value = 8008.4688
$ws.Range("L32").Value = 12534.546
$ws.Range("M32").Value = -7721.4688
$ws.Range("N32").Value = -13108.546
$ws.Range("H132").Value = 7145301
$ws.Range("I132").Value = 13160089
$ws.Range("K132").Value = 39480267
$ws.Range("M132").Value = -39477737

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 18600
$ws.Range("J39").Value = 18600
$ws.Range("L39").Value = 18600
$ws.Range("N39").Value = -19378
$ws.Range("H56").Value = 49582.5
$ws.Range("J56").Value = 49582.5
$ws.Range("L56").Value = 49582.5
$ws.Range("N56").Value = -51060.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 32515.75
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 32515.75
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 32515.75
$ws.Range("N26").Value = -33089.75
$ws.Range("M26").ClearContents()
$ws.Range("H35").Value = 615
$ws.Range("I35").Value = 529.44446
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 529.44446
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = -235.44446
$ws.Range("N35").Value = -1588
$ws.Range("H44").Value = 170000
$ws.Range("I44").Value = 170000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 170000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -169558
$ws.Range("N44").ClearContents()
$ws.Range("H53").Value = 23300
$ws.Range("J53").Value = 23300
$ws.Range("L53").Value = 23300
$ws.Range("N53").Value = -24514
$ws.Range("H54").Value = 29955.2
$ws.Range("J54").Value = 29955.2
$ws.Range("L54").Value = 29955.2
$ws.Range("N54").Value = -31271.2
$ws.Range("H56").Value = 33000
$ws.Range("J56").Value = 33000
$ws.Range("L56").Value = 33000
$ws.Range("N56").Value = -34690
$ws.Range("H58").Value = 1414.6531
$ws.Range("I58").Value = 747.25
$ws.Range("J58").Value = 2304.524
$ws.Range("K58").Value = 747.25
$ws.Range("L58").Value = 2304.524
$ws.Range("M58").Value = -544.25
$ws.Range("N58").Value = -2710.524
$ws.Range("H132").Value = 8773668
$ws.Range("I132").Value = 11906138
$ws.Range("J132").Value = 2754.5334
$ws.Range("K132").Value = 35718414
$ws.Range("L132").Value = 8263.600199999999
$ws.Range("M132").Value = -35715884
$ws.Range("N132").Value = -13323.6002
$ws.Range("H134").Value = 1512
$ws.Range("I134").Value = 1652.1904
$ws.Range("J134").Value = 1091.4286
$ws.Range("K134").Value = 4956.5712
$ws.Range("L134").Value = 3274.2858
$ws.Range("M134").Value = -2421.5712
$ws.Range("N134").Value = -8344.2858
$ws.Range("H136").Value = 1414.6531
$ws.Range("I136").Value = 747.25
$ws.Range("J136").Value = 2304.524
$ws.Range("K136").Value = 2241.75
$ws.Range("L136").Value = 6913.572
$ws.Range("M136").Value = 308.25
$ws.Range("N136").Value = -12013.572

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5749.75
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 5749.75
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 17249.25
$ws.Range("N42").Value = -18317.25
$ws.Range("M42").ClearContents()
$ws.Range("H98").Value = 1075.5555
$ws.Range("J98").Value = 1463.3334
$ws.Range("L98").Value = 4390.0002
$ws.Range("N98").Value = -7386.0002
$ws.Range("H99").Value = 2102.5
$ws.Range("I99").Value = 512.5
$ws.Range("K99").Value = 1537.5
$ws.Range("M99").Value = 708.5
$ws.Range("H101").Value = 11466.667
$ws.Range("J101").Value = 11466.667
$ws.Range("L101").Value = 34400.001
$ws.Range("N101").Value = -39268.001
$ws.Range("H131").Value = 833.97
$ws.Range("J131").Value = 859.10754
$ws.Range("L131").Value = 2577.32262
$ws.Range("N131").Value = -12657.32262

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3362.7754
$ws.Range("J126").Value = 4529.7617
$ws.Range("L126").Value = 13589.2851
$ws.Range("N126").Value = -18529.2851

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1423
$ws.Range("I46").Value = 926.6667
$ws.Range("J46").Value = 1635.7142
$ws.Range("K46").Value = 926.6667
$ws.Range("L46").Value = 1635.7142
$ws.Range("M46").Value = -738.6667
$ws.Range("N46").Value = -2011.7142
$ws.Range("H68").Value = 2377.7778
$ws.Range("I68").Value = 2600
$ws.Range("J68").Value = 2266.6667
$ws.Range("K68").Value = 2600
$ws.Range("L68").Value = 2266.6667
$ws.Range("M68").Value = -1851
$ws.Range("N68").Value = -3764.6667
$ws.Range("H71").Value = 2377.7778
$ws.Range("I71").Value = 2600
$ws.Range("J71").Value = 2266.6667
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 11333.3335
$ws.Range("M71").Value = -9256
$ws.Range("N71").Value = -18821.3335
$ws.Range("H94").Value = 50631.668
$ws.Range("J94").Value = 50631.668
$ws.Range("L94").Value = 50631.668
$ws.Range("N94").Value = -51983.668
$ws.Range("H132").Value = 9440620
$ws.Range("I132").Value = 3842.0571
$ws.Range("J132").Value = 27789912
$ws.Range("K132").Value = 11526.1713
$ws.Range("L132").Value = 83369736
$ws.Range("M132").Value = -8996.1713
$ws.Range("N132").Value = -83374796

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 14810.909
$ws.Range("J33").Value = 16192
$ws.Range("L33").Value = 16192
$ws.Range("N33").Value = -16692
$ws.Range("H36").Value = 14810.909
$ws.Range("J36").Value = 16192
$ws.Range("L36").Value = 16192
$ws.Range("N36").Value = -16692
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H42").Value = 38900
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 38900
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 38900
$ws.Range("N42").Value = -39656
$ws.Range("M42").ClearContents()
$ws.Range("H62").Value = 10230.23
$ws.Range("I62").Value = 4931.25
$ws.Range("J62").Value = 18708.6
$ws.Range("K62").Value = 4931.25
$ws.Range("L62").Value = 18708.6
$ws.Range("M62").Value = -4307.25
$ws.Range("N62").Value = -19956.6
$ws.Range("H65").Value = 10230.23
$ws.Range("I65").Value = 4931.25
$ws.Range("J65").Value = 18708.6
$ws.Range("K65").Value = 24656.25
$ws.Range("L65").Value = 93543
$ws.Range("M65").Value = -21536.25
$ws.Range("N65").Value = -99783
$ws.Range("H105").Value = 19500
$ws.Range("J105").Value = 19500
$ws.Range("L105").Value = 19500
$ws.Range("N105").Value = -26488
$ws.Range("H136").Value = 986.3
$ws.Range("I136").Value = 849.18604
$ws.Range("J136").Value = 1828.5714
$ws.Range("K136").Value = 2547.55812
$ws.Range("L136").Value = 5485.7142
$ws.Range("M136").Value = 2.441879999999855
$ws.Range("N136").Value = -10585.7142
